$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update per-country statistics (new data for this refresh) ---------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes
$updates = @{
    "Israel"    = @(4347, 100, 134, 4198, 80, 0, 15)
    "Chequia"   = @(2837, 20, 11, 2809, 45, 1, 17)
    "Estonia"   = @(715, 36, 20, 692, 10, 0, 3)
    "Marruecos" = @(516, 37, 13, 476, 1, 1, 27)
    "Barein"    = @(515, 16, 279, 232, 2, 0, 4)
    "Letonia"   = @(376, 29, 1, 375, 3, 0, 0)
    "Georgia"   = @(98, 7, 18, 80, 1, 0, 0)
}

$lastRow = $ws.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count

foreach ($country in $updates.Keys) {
    $found = $ws.Range("A4:A205").Find($country)
    if ($found -ne $null) {
        $r = $found.Row
        $vals = $updates[$country]
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
        $ws.Cells.Item($r, 5).Value = $vals[3]
        $ws.Cells.Item($r, 6).Value = $vals[4]
        $ws.Cells.Item($r, 7).Value = $vals[5]
        $ws.Cells.Item($r, 8).Value = $vals[6]
    }
}

# --- 2. Re-sort the country table by "Casos totales" (column B), descending
$dataRange = $ws.Range("A4:H205")
$sortKey = $ws.Range("B4:B205")
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1)

# --- 3. Update the "last updated" timestamp banner --------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 09:50"
